# Weekly update: insert a new week's worth of data (2 rows: "Primera" and
# "Segunda" quality) at the top of the Betarraga price table, pushing all
# existing rows down by two (matches "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the first data row of this block (310),
# shifting rows 310:389 down to 312:391 (and growing the used range to R391).
$ws.Rows("310:311").Insert()

# New row 310 - "Primera" quality for the new week.
$ws.Cells.Item(310, 1).Value = 8
$ws.Cells.Item(310, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(310, 3).Value = "Coquimbo"
$ws.Cells.Item(310, 4).Value = 44943
$ws.Cells.Item(310, 5).Value = 4
$ws.Cells.Item(310, 6).Value = 100114014
$ws.Cells.Item(310, 7).Value = "Betarraga"
$ws.Cells.Item(310, 8).Value = "Sin especificar"
$ws.Cells.Item(310, 9).Value = "Primera"
$ws.Cells.Item(310, 10).Value = 2000
$ws.Cells.Item(310, 11).Value = 500
$ws.Cells.Item(310, 12).Value = 600
$ws.Cells.Item(310, 13).Value = 550
$ws.Cells.Item(310, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(310, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(310, 16).Value = 183
$ws.Cells.Item(310, 17).Value = 3
$ws.Cells.Item(310, 18).Value = "Hortaliza"

# New row 311 - "Segunda" quality for the new week.
$ws.Cells.Item(311, 1).Value = 8
$ws.Cells.Item(311, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(311, 3).Value = "Coquimbo"
$ws.Cells.Item(311, 4).Value = 44943
$ws.Cells.Item(311, 5).Value = 4
$ws.Cells.Item(311, 6).Value = 100114014
$ws.Cells.Item(311, 7).Value = "Betarraga"
$ws.Cells.Item(311, 8).Value = "Sin especificar"
$ws.Cells.Item(311, 9).Value = "Segunda"
$ws.Cells.Item(311, 10).Value = 1560
$ws.Cells.Item(311, 11).Value = 400
$ws.Cells.Item(311, 12).Value = 450
$ws.Cells.Item(311, 13).Value = 425
$ws.Cells.Item(311, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(311, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(311, 16).Value = 142
$ws.Cells.Item(311, 17).Value = 3
$ws.Cells.Item(311, 18).Value = "Hortaliza"
